# Add 2022-Q4 data
# 1) Insert a new row at the top of the "总计" (totals) sheet's data for 2022-Q4.
# 2) Insert a brand-new worksheet named "2022-Q4" right after "总计" (i.e. before
#    what is currently "2022-Q3"), populated with the new quarter's fund data.
# Every other existing sheet keeps its own content and simply shifts down by one
# tab position, which Excel does automatically when a sheet is inserted/moved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for 2022-Q4 and push the rest down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Copy the formatting from row 3 (an existing data row) onto the new row 2
# first, so the styling (centered index column, plain data columns) matches
# before the values are written in.
$total.Range("A3:D3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.58

# Renumber the index column (A) for the rows that used to start at 0.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right before "2022-Q3".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

$fundRows = @(
  @("000593", "易方达标普全球高端消费品指数增强（QDII）美元现汇"),
  @("005676", "易方达标普全球高端消费品指数增强C（QDII）人民币"),
  @("118002", "易方达标普全球高端消费品指数增强A（QDII）人民币")
)

$r = 2
foreach ($fr in $fundRows) {
    $q4.Range("B$r").NumberFormat = "@"
    $q4.Range("B$r").Value = $fr[0]
    $q4.Range("C$r").Value = $fr[1]
    $q4.Range("D$r").NumberFormat = "@"
    $q4.Range("D$r").Value = "2.30"
    $q4.Range("E$r").NumberFormat = "@"
    $q4.Range("E$r").Value = "93.71"
    $q4.Range("F$r").NumberFormat = "@"
    $q4.Range("F$r").Value = "8.44"
    $q4.Range("G$r").NumberFormat = "@"
    $q4.Range("G$r").Value = "0.1941"
    $q4.Range("H$r").Value = 4
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Restore the originally-active tab (last sheet, "2021-Q2").
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
